# Updates the cryptocurrency price/volume table in Sheet1 (columns D, E)
# to the latest scraped snapshot. A couple of coins (rows 15/16) also
# swapped rank position, so their Name/Link/Price/Volume cells moved.
#
# Price/Volume cells are stored as literal text in the workbook (e.g.
# "201.60", "0.999", "  +1.33%  ") rather than numbers, so that trailing
# zeros / leading "+" signs / padding spaces render exactly as scraped.
# Assigning a numeric-looking string straight to Range.Value lets Excel
# auto-coerce it to a real number (dropping formatting like the trailing
# zero in "201.60"), so for any new value that looks numeric we briefly
# mark the cell as Text ("@") before writing it, then restore the
# original "Normal" style so no extra formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '76.483.37'
$ws.Range("E2").Value = '  +1.33%  '

# Row 3
$ws.Range("D3").Value = '2.941.56'
$ws.Range("E3").Value = '  +4.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '598.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("E8").Value = '  +1.24%  '

# Row 9
$ws.Range("E9").Value = '  +2.63%  '

# Row 10
$ws.Range("D10").Value = '2.941.27'
$ws.Range("E10").Value = '  +4.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.443'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.31%  '

# Row 12
$ws.Range("E12").Value = '  +0.86%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.23%  '

# Row 14
$ws.Range("D14").Value = '3.480.86'
$ws.Range("E14").Value = '  +4.29%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '76.345.26'
$ws.Range("E15").Value = '  +1.41%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.35%  '

# Row 17
$ws.Range("E17").Value = '  +1.20%  '

# Row 18
$ws.Range("D18").Value = '2.930.82'
$ws.Range("E18").Value = '  +3.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.16%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '

# Row 22
$ws.Range("E22").Value = '  -0.33%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '

# Row 26
$ws.Range("D26").Value = '3.085.93'
$ws.Range("E26").Value = '  +4.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.83%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000109'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.44%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("E31").Value = '  -1.34%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.92%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '499.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.10%  '

# Row 34
$ws.Range("E34").Value = '  +1.70%  '

# Row 35
$ws.Range("E35").Value = '  +0.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.393'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.60%  '

# Row 39
$ws.Range("E39").Value = '  +23.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.42%  '

# Row 41
$ws.Range("E41").Value = '  -5.26%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '178.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.12%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.03%  '

# Row 45
$ws.Range("E45").Value = '  -0.72%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.65%  '

# Row 47
$ws.Range("E47").Value = '  -0.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.87%  '

# Row 49
$ws.Range("E49").Value = '  +2.57%  '

# Row 50
$ws.Range("E50").Value = '  +3.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.60%  '
